# Round 3 final setup
# Updates Round 1 / Round 2 scores for the top of the Intermediate leaderboard
# and normalises the capitalisation of a few player names.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Score corrections (Round 1 = column C, Round 2 = column D) ---
$ws.Range("D2").Value = 108000

$ws.Range("C3").Value = 11600
$ws.Range("D3").Value = 42800

$ws.Range("C4").Value = 80100
$ws.Range("C5").Value = 80000
$ws.Range("C6").Value = 66000
$ws.Range("C7").Value = 60800
$ws.Range("C8").Value = 50800
$ws.Range("C9").Value = 41000

$ws.Range("D10").Value = 33200

$ws.Range("C11").Value = 23600
$ws.Range("C12").Value = 3900
$ws.Range("C13").Value = -16000

# --- Player name capitalisation fixes ---
$ws.Range("B18").Value = "Jonah wentzel"
$ws.Range("B25").Value = "brady surya sie"
$ws.Range("B34").Value = "Ian carroll"

# --- Restore the last-used selection ---
$ws.Range("A2:A35").Select()
